$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 3 de Julio de 2020 a las 03:08"

# Row 4: Estados Unidos - update numeric stats
$ws.Range("B4").Value = 2837189
$ws.Range("C4").Value = 57236
$ws.Range("D4").Value = 1191091
$ws.Range("E4").Value = 1514613
$ws.Range("G4").Value = 687
$ws.Range("H4").Value = 131485

# Row 23: Canada - update numeric stats
$ws.Range("B23").Value = 104772
$ws.Range("C23").Value = 501
$ws.Range("E23").Value = 27783

# Rows 82/83: Venezuela overtakes Kirguistan in ranking
$ws.Range("A82").Value = "Venezuela"
$ws.Range("B82").Value = 6273
$ws.Range("C82").Value = 211
$ws.Range("D82").Value = 2100
$ws.Range("E82").Value = 4116
$ws.Range("G82").Value = 3
$ws.Range("H82").Value = 57

$ws.Range("A83").Value = "Kirguistan"
$ws.Range("B83").Value = 6261
$ws.Range("C83").Value = 526
$ws.Range("D83").Value = 2530
$ws.Range("E83").Value = 3665
$ws.Range("G83").Value = 4
$ws.Range("H83").Value = 66

# Row 96: Costa Rica - update numeric stats
$ws.Range("E96").Value = 2416
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 18

# Rows 203-206: Laos/Santa Lucia and Dominica/Fiyi swap (tied totals, cosmetic reorder)
$ws.Range("A203").Value = "Laos"
$ws.Range("A204").Value = "Santa Lucia"
$ws.Range("A205").Value = "Dominica"
$ws.Range("A206").Value = "Fiyi"

# Rows 209-210: Islas Malvinas/Groenlandia swap (tied totals, cosmetic reorder)
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"
